# The commit swaps the two theme parts in the package: the theme that is
# actually wired up to the slide master / presentation (ppt/theme/theme2.xml,
# currently the "Integral" / "Red Violet" palette) ends up holding the plain
# "Office Theme" / "Office" palette that used to live in the unused
# ppt/theme/theme1.xml part (which is only referenced by the notes master
# and isn't reachable through the Slide/Master/Design object model).
#
# The PowerPoint object model lets us reach the live theme through
# Master.Theme.ThemeColorScheme (12 color slots, same order as a:clrScheme:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink). Setting .RGB on each slot
# rewrites ppt/theme/theme2.xml's <a:srgbClr val="..."/> entries in place,
# which is exactly the part of the diff we can reproduce.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

function Set-ThemeColor($scheme, $index, $r, $g, $b) {
    $oleColor = ($b * 65536) + ($g * 256) + $r
    $scheme.Item($index).RGB = $oleColor
}

# Target palette = the "Office Theme" colors (previously sitting unused in
# ppt/theme/theme1.xml), applied in clrScheme slot order.
Set-ThemeColor $colorScheme 1  0   0   0     # dk1      000000
Set-ThemeColor $colorScheme 2  255 255 255   # lt1      FFFFFF
Set-ThemeColor $colorScheme 3  0x44 0x54 0x6A # dk2     44546A
Set-ThemeColor $colorScheme 4  0xE7 0xE6 0xE6 # lt2     E7E6E6
Set-ThemeColor $colorScheme 5  0x5B 0x9B 0xD5 # accent1 5B9BD5
Set-ThemeColor $colorScheme 6  0xED 0x7D 0x31 # accent2 ED7D31
Set-ThemeColor $colorScheme 7  0xA5 0xA5 0xA5 # accent3 A5A5A5
Set-ThemeColor $colorScheme 8  0xFF 0xC0 0x00 # accent4 FFC000
Set-ThemeColor $colorScheme 9  0x44 0x72 0xC4 # accent5 4472C4
Set-ThemeColor $colorScheme 10 0x70 0xAD 0x47 # accent6 70AD47
Set-ThemeColor $colorScheme 11 0x05 0x63 0xC1 # hlink   0563C1
Set-ThemeColor $colorScheme 12 0x95 0x4F 0x72 # folHlink 954F72
